# Apply crypto symbol-list update
# (commit: Updated symbol list on Sun Dec 25 06:13:10 UTC 2022 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores numeric-looking values
    # (prices, hour counters) as literal text, same as the source sheet.
    $ws.Range($cellRef).Value = "'" + $text
}

# Row 2
Set-TextCell "D2" "245.05"
Set-TextCell "G2" "6"
# Row 3
Set-TextCell "D3" "23.02"
Set-TextCell "G3" "6"
# Row 4
Set-TextCell "D4" "5.408"
Set-TextCell "G4" "6"
# Row 5
Set-TextCell "D5" "0.06020"
Set-TextCell "G5" "6"
# Row 6
Set-TextCell "G6" "6"
# Row 7
Set-TextCell "D7" "0.8102"
Set-TextCell "G7" "6"
# Row 8
Set-TextCell "D8" "0.9295"
Set-TextCell "G8" "6"
# Row 9
Set-TextCell "D9" "0.1422"
Set-TextCell "G9" "6"
# Row 10
Set-TextCell "D10" "0.07432"
Set-TextCell "G10" "6"
# Row 11
Set-TextCell "D11" "0.03381"
Set-TextCell "G11" "6"
# Row 12
Set-TextCell "G12" "6"
# Row 13
Set-TextCell "D13" "0.09358"
Set-TextCell "G13" "6"
# Row 14
Set-TextCell "D14" "3.936"
Set-TextCell "G14" "6"
# Row 15
Set-TextCell "D15" "0.001590"
Set-TextCell "G15" "6"
# Row 16
Set-TextCell "D16" "0.04835"
Set-TextCell "G16" "6"
# Row 17
Set-TextCell "G17" "6"
# Row 18
Set-TextCell "D18" "0.005379"
Set-TextCell "G18" "6"
# Row 19
Set-TextCell "D19" "0.004154"
Set-TextCell "G19" "6"
# Row 20
Set-TextCell "D20" "0.0009840"
Set-TextCell "G20" "6"
# Row 21
Set-TextCell "D21" "0.00008703"
Set-TextCell "G21" "6"
# Row 22
Set-TextCell "D22" "3.653"
Set-TextCell "G22" "6"
# Row 23
Set-TextCell "D23" "6.442"
Set-TextCell "G23" "6"
# Row 24
Set-TextCell "D24" "2.186"
Set-TextCell "G24" "6"
# Row 25
Set-TextCell "G25" "6"
# Row 26
Set-TextCell "D26" "0.1294"
Set-TextCell "G26" "6"
# Row 27
Set-TextCell "G27" "6"
# Row 28
Set-TextCell "G28" "6"
# Row 29
Set-TextCell "G29" "6"
# Row 30
Set-TextCell "G30" "6"
# Row 31
Set-TextCell "G31" "6"
# Row 32
Set-TextCell "G32" "6"
# Row 33
Set-TextCell "G33" "6"
# Row 34
Set-TextCell "G34" "6"
# Row 35
Set-TextCell "G35" "6"
# Row 36
Set-TextCell "G36" "6"
# Row 37
Set-TextCell "G37" "6"
# Row 38
Set-TextCell "G38" "6"
# Row 39
Set-TextCell "G39" "6"
# Row 40
Set-TextCell "D40" "0.03972"
Set-TextCell "G40" "6"
# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D41" "0.006372"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextCell "G41" "6"
# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
Set-TextCell "G42" "6"
# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D43" "0.002711"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextCell "G43" "6"
# Row 44
Set-TextCell "D44" "0.006588"
Set-TextCell "G44" "6"
# Row 45
Set-TextCell "D45" "0.00005213"
Set-TextCell "G45" "6"
# Row 46
Set-TextCell "G46" "6"
# Row 47
Set-TextCell "G47" "6"
# Row 48
Set-TextCell "D48" "0.8203"
Set-TextCell "G48" "6"
# Row 49
Set-TextCell "D49" "0.002299"
Set-TextCell "G49" "6"
# Row 50
Set-TextCell "G50" "6"
# Row 51
Set-TextCell "G51" "6"
